$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new weekly observation rows at row 208; everything currently at
# row 208 and below (through 231) shifts down to 210..233, which matches the
# newer data shown further down in the sheet (rows 210..233 == old 208..231).
$ws.Rows("208:209").Insert()

# New row 208: Chirimoya, Especial quality, week of 45212
$ws.Cells.Item(208, 1).Value  = 10
$ws.Cells.Item(208, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(208, 3).Value  = "La Araucanía"
$ws.Cells.Item(208, 4).Value  = 45212
$ws.Cells.Item(208, 5).Value  = 9
$ws.Cells.Item(208, 6).Value  = "Fruta"
$ws.Cells.Item(208, 7).Value  = 100107
$ws.Cells.Item(208, 8).Value  = "Otros"
$ws.Cells.Item(208, 9).Value  = 100107002
$ws.Cells.Item(208, 10).Value = "Chirimoya"
$ws.Cells.Item(208, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(208, 12).Value = "Especial"
$ws.Cells.Item(208, 13).Value = 70
$ws.Cells.Item(208, 14).Value = 3000
$ws.Cells.Item(208, 15).Value = 3000
$ws.Cells.Item(208, 16).Value = 3000
$ws.Cells.Item(208, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(208, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(208, 19).Value = 3000
$ws.Cells.Item(208, 20).Value = 1

# New row 209: Chirimoya, Primera quality, week of 45212
$ws.Cells.Item(209, 1).Value  = 10
$ws.Cells.Item(209, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(209, 3).Value  = "La Araucanía"
$ws.Cells.Item(209, 4).Value  = 45212
$ws.Cells.Item(209, 5).Value  = 9
$ws.Cells.Item(209, 6).Value  = "Fruta"
$ws.Cells.Item(209, 7).Value  = 100107
$ws.Cells.Item(209, 8).Value  = "Otros"
$ws.Cells.Item(209, 9).Value  = 100107002
$ws.Cells.Item(209, 10).Value = "Chirimoya"
$ws.Cells.Item(209, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(209, 12).Value = "Primera"
$ws.Cells.Item(209, 13).Value = 170
$ws.Cells.Item(209, 14).Value = 2500
$ws.Cells.Item(209, 15).Value = 2600
$ws.Cells.Item(209, 16).Value = 2547
$ws.Cells.Item(209, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(209, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(209, 19).Value = 2547
$ws.Cells.Item(209, 20).Value = 1
